$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume-change (column E) values.
# NumberFormat is forced to Text ("@") before assignment so Excel does not
# auto-convert numeric-looking strings (e.g. "41.200.75", "2.10", "0.0998")
# into numbers and strip formatting such as trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.200.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.222.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.51%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.55"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.90%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.34"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.51%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.556.52"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.80"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -9.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -8.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.227.18"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.230.90"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0957"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.61"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -8.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +12.16%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.81%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.07"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.56"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.87%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.97%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.38%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -10.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.99"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +12.28%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.86"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -11.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.84"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.02"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.19%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.77"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.85"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +10.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0998"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.97%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.53"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.72%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.24%  "
